# The commit swaps the two theme parts of the deck: the theme that is
# actually "live" (applied to the slide master / whole deck) switches from
# the "Integral" / Red Violet palette back to the stock Office Theme
# palette (the other theme part -- only ever linked from the notes master
# -- picks up the Integral colours, but that part is not reachable through
# the PowerPoint object model, so here we reproduce the visible half of the
# swap: the colours that are actually rendered throughout the presentation).

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colour scheme, in ThemeColorScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeTheme = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeTheme[$i - 1]
}

Write-Output ("Updated " + $colors.Count + " theme colors on '" + $master.Name + "'")
